$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("필수과목목록")

# Fill the previously-empty "대체교과목번호" cells (F2:F15) with a single
# space character, matching the placeholder value already used elsewhere
# in the sheet (e.g. column B).
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 6).Value = " "
}

# Move the active selection to F18, matching the cell last touched.
$ws.Range("F18").Select()
